# QA/Property_Booking_Bug_list.xlsx - "updated QA list"
#
# Bug #11 (row 13) used to be tracked separately from bug #12 (row 14):
#   H13 = "Resolved e,f"   (the "send mail form" items e & f)
#   H14 = "Resolved c,d"   (the "send mail form" items c & d)
# The QA list was updated to roll all four fixes into the row-13 status
# cell and clear the note back down to the bare status on row 14, and a
# "Resolved" status was filled in for bug #17 (row 19), which had been
# left blank. Several status cells also picked up wrap-text formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# H12 ("Resolved") gains wrap-text formatting; value is unchanged.
$ws.Range("H12").Value = "Resolved"
$ws.Range("H12").WrapText = $true

# H13: "Resolved e,f" -> "Resolved c,d,e,f", plus wrap-text formatting.
$ws.Range("H13").Value = "Resolved c,d,e,f"
$ws.Range("H13").WrapText = $true

# H14: trimmed back down to "Resolved c,d".
$ws.Range("H14").Value = "Resolved c,d"

# H19 was empty; bug #17 is now marked Resolved.
$ws.Range("H19").Value = "Resolved"

# Scroll the sheet so row 11 is at the top and H14 is the active selection.
$ws.Activate()
$excel.Goto($ws.Range("A11"), $true)
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H14").Select()
